# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1884
    "F5"  = 195
    "F6"  = 754
    "F7"  = 305
    "F8"  = 361
    "F9"  = 4498
    "F10" = 23
    "F11" = 353
    "F12" = 1274
    "F13" = 530
    "F15" = 859
    "F17" = 473
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
